$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.982.78'
$ws.Cells.Item(2, 5).Value = '  -0.70%  '

$ws.Cells.Item(3, 4).Value = '2.642.83'
$ws.Cells.Item(3, 5).Value = '  +0.36%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = '581.35'
$ws.Cells.Item(5, 5).Value = '  +0.06%  '

$ws.Cells.Item(6, 4).Value = '156.64'
$ws.Cells.Item(6, 5).Value = '  +0.15%  '

$ws.Cells.Item(7, 4).Value = '0.626'
$ws.Cells.Item(7, 5).Value = '  -3.05%  '

$ws.Cells.Item(8, 5).Value = '  +0.01%  '

$ws.Cells.Item(9, 4).Value = '2.640.98'
$ws.Cells.Item(9, 5).Value = '  +0.36%  '

$ws.Cells.Item(10, 4).Value = '0.118'
$ws.Cells.Item(10, 5).Value = '  -3.39%  '

$ws.Cells.Item(11, 5).Value = '  +0.01%  '

$ws.Cells.Item(12, 5).Value = '  -1.07%  '

$ws.Cells.Item(13, 4).Value = '0.156'
$ws.Cells.Item(13, 5).Value = '  +0.84%  '

$ws.Cells.Item(14, 4).Value = '28.65'
$ws.Cells.Item(14, 5).Value = '  -0.62%  '

$ws.Cells.Item(15, 4).Value = '3.119.30'
$ws.Cells.Item(15, 5).Value = '  +0.23%  '

$ws.Cells.Item(16, 5).Value = '  -0.97%  '

$ws.Cells.Item(17, 4).Value = '63.898.87'
$ws.Cells.Item(17, 5).Value = '  -0.53%  '

$ws.Cells.Item(18, 4).Value = '2.639.63'
$ws.Cells.Item(18, 5).Value = '  -0.38%  '

$ws.Cells.Item(19, 4).Value = '12.19'
$ws.Cells.Item(19, 5).Value = '  -0.14%  '

$ws.Cells.Item(20, 4).Value = '7.84'
$ws.Cells.Item(20, 5).Value = '  +3.48%  '

$ws.Cells.Item(21, 5).Value = '  -3.35%  '

$ws.Cells.Item(22, 4).Value = '345.80'
$ws.Cells.Item(22, 5).Value = '  -0.52%  '

$ws.Cells.Item(23, 5).Value = '  +0.17%  '

$ws.Cells.Item(24, 4).Value = '68.20'
$ws.Cells.Item(24, 5).Value = '  +0.16%  '

$ws.Cells.Item(25, 4).Value = '1.86'
$ws.Cells.Item(25, 5).Value = '  +4.66%  '

$ws.Cells.Item(26, 4).Value = '0.0000112'
$ws.Cells.Item(26, 5).Value = '  +1.90%  '

$ws.Cells.Item(27, 4).Value = '9.39'
$ws.Cells.Item(27, 5).Value = '  -0.37%  '

$ws.Cells.Item(28, 4).Value = '1.65'
$ws.Cells.Item(28, 5).Value = '  +3.57%  '

$ws.Cells.Item(29, 4).Value = '589.04'
$ws.Cells.Item(29, 5).Value = '  -0.34%  '

$ws.Cells.Item(30, 4).Value = '8.21'
$ws.Cells.Item(30, 5).Value = '  +2.35%  '

$ws.Cells.Item(31, 5).Value = '  -0.03%  '

$ws.Cells.Item(32, 4).Value = '0.999'
$ws.Cells.Item(32, 5).Value = '  -0.21%  '

$ws.Cells.Item(33, 5).Value = '  -0.58%  '

$ws.Cells.Item(34, 4).Value = '1.74'
$ws.Cells.Item(34, 5).Value = '  +0.57%  '

$ws.Cells.Item(35, 5).Value = '  +0.23%  '

$ws.Cells.Item(36, 4).Value = '5.50'
$ws.Cells.Item(36, 5).Value = '  +3.38%  '

$ws.Cells.Item(37, 5).Value = '  -2.30%  '

$ws.Cells.Item(38, 4).Value = '19.79'
$ws.Cells.Item(38, 5).Value = '  -1.08%  '

$ws.Cells.Item(39, 5).Value = '  -0.04%  '

$ws.Cells.Item(40, 5).Value = '  -0.45%  '

$ws.Cells.Item(41, 4).Value = '151.64'
$ws.Cells.Item(41, 5).Value = '  +0.55%  '

$ws.Cells.Item(42, 4).Value = '2.59'
$ws.Cells.Item(42, 5).Value = '  +9.28%  '

$ws.Cells.Item(43, 5).Value = '  +0.00%  '

$ws.Cells.Item(44, 4).Value = '41.99'
$ws.Cells.Item(44, 5).Value = '  +0.02%  '

$ws.Cells.Item(45, 4).Value = '163.90'
$ws.Cells.Item(45, 5).Value = '  +3.01%  '

$ws.Cells.Item(46, 4).Value = '24.45'
$ws.Cells.Item(46, 5).Value = '  +4.10%  '

$ws.Cells.Item(47, 4).Value = '3.91'
$ws.Cells.Item(47, 5).Value = '  -2.62%  '

$ws.Cells.Item(48, 4).Value = '0.0590'
$ws.Cells.Item(48, 5).Value = '  -2.09%  '

$ws.Cells.Item(49, 5).Value = '  +0.10%  '

$ws.Cells.Item(50, 5).Value = '  -2.48%  '

$ws.Cells.Item(51, 5).Value = '  -1.88%  '
